# Slide 21 ("Find the syndrome..." / Hamming-code example): the
# "Message bit" callout text box was updated from a 3-bit message to a
# 4-bit message, and (since the box wraps text at a fixed width rather
# than auto-sizing horizontally) the author widened the box so the new,
# longer line still fits on one line.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(21)
$sh = $s.Shapes.Item("TextBox 7")

# Update the caption text (keeps existing run formatting - font/size).
$sh.TextFrame.TextRange.Text = "Message bit [1 0 1 0]"

# Widen the text box to match the author's manual resize
# (2253343 EMU -> 2616015 EMU; height/position are unchanged).
$sh.Width = 205.98547244094487
